$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift accelerometer/gyroscope readings down by two rows and splice in two new
# leading samples (falling class), then drop the now-superfluous last row.
$arr = New-Object 'double[,]' 20,6
$arr[0,0] = -2.092850303649902
$arr[0,1] = 0.9900987625122062
$arr[0,2] = 1.828120517730715
$arr[0,3] = 0.0937678143382072
$arr[0,4] = 0.0774271711707115
$arr[0,5] = -0.1965458989143371
$arr[1,0] = -2.105730056762695
$arr[1,1] = 0.8735208511352541
$arr[1,2] = 2.453470587730407
$arr[1,3] = -0.1837177276611328
$arr[1,4] = 0.2562579810619354
$arr[1,5] = 0.0125227374956011
$arr[2,0] = -2.657281875610353
$arr[2,1] = 0.2476238250732379
$arr[2,2] = 2.285648679733275
$arr[2,3] = -0.0036651915870606
$arr[2,4] = 0.1434006094932556
$arr[2,5] = -0.0604756586253643
$arr[3,0] = -4.063823699951181
$arr[3,1] = -1.187473011016849
$arr[3,2] = 2.541466045379643
$arr[3,3] = 0.2081523388624191
$arr[3,4] = 0.1316414624452591
$arr[3,5] = 0.1007927656173706
$arr[4,0] = -5.616066837310788
$arr[4,1] = -2.403632545471194
$arr[4,2] = 3.652213478088379
$arr[4,3] = -0.4978551864624023
$arr[4,4] = -0.2622139155864715
$arr[4,5] = -1.034500241279602
$arr[5,0] = -3.539181804656983
$arr[5,1] = -3.395954704284656
$arr[5,2] = 3.008686828613268
$arr[5,3] = -0.328340083360672
$arr[5,4] = -0.4218024611473083
$arr[5,5] = -0.5737552046775818
$arr[6,0] = -3.633686542510985
$arr[6,1] = 3.542618751525879
$arr[6,2] = -4.68873119354248
$arr[6,3] = -2.957962274551392
$arr[6,4] = -0.2176207453012466
$arr[6,5] = -2.633592844009399
$arr[7,0] = 53.59859118461612
$arr[7,1] = -4.54958009719849
$arr[7,2] = -34.16661596298219
$arr[7,3] = 1.301142930984497
$arr[7,4] = -3.674507141113281
$arr[7,5] = 2.831360340118408
$arr[8,0] = 10.96323089599605
$arr[8,1] = -0.1360907554626412
$arr[8,2] = -11.01923332214353
$arr[8,3] = 0.7254024744033813
$arr[8,4] = 0.9521862268447876
$arr[8,5] = -1.527926683425903
$arr[9,0] = -1.288498878479004
$arr[9,1] = 0.4897777378559073
$arr[9,2] = -1.382943773269636
$arr[9,3] = -0.2622139155864715
$arr[9,4] = 1.557553648948669
$arr[9,5] = -0.6866125464439392
$arr[10,0] = 0.2403127670288223
$arr[10,1] = 0.4333343148231529
$arr[10,2] = 1.078743743896495
$arr[10,3] = 1.219134330749511
$arr[10,4] = -0.9372199773788452
$arr[10,5] = 2.843730449676514
$arr[11,0] = 0.6345248222351046
$arr[11,1] = 0.4263583719730371
$arr[11,2] = 1.455561161041258
$arr[11,3] = 0.8478809595108032
$arr[11,4] = -0.9918924570083618
$arr[11,5] = -0.8237518072128296
$arr[12,0] = -0.09080390930175844
$arr[12,1] = 2.436175584793104
$arr[12,2] = 1.085069131851197
$arr[12,3] = -0.0302378293126821
$arr[12,4] = -0.6291912198066711
$arr[12,5] = -0.0180205255746841
$arr[13,0] = -1.009081411361702
$arr[13,1] = 3.721519541740399
$arr[13,2] = 1.14239126443863
$arr[13,3] = -0.0595593601465225
$arr[13,4] = -0.2151772826910019
$arr[13,5] = 0.1901318132877349
$arr[14,0] = -2.444419670104975
$arr[14,1] = -0.2522135257720842
$arr[14,2] = 1.136495351791381
$arr[14,3] = 0.0241291765123605
$arr[14,4] = 0.3474296033382416
$arr[14,5] = 0.1518000066280365
$arr[15,0] = -1.008091449737549
$arr[15,1] = 2.286790394783017
$arr[15,2] = 0.9021725535392736
$arr[15,3] = -0.0592539273202419
$arr[15,4] = 0.4436408877372741
$arr[15,5] = 0.2559525370597839
$arr[16,0] = -1.014060974121094
$arr[16,1] = 1.50408124923706
$arr[16,2] = 0.3019133806228637
$arr[16,3] = 0.057115901261568
$arr[16,4] = 0.0777326002717018
$arr[16,5] = -0.1259909570217132
$arr[17,0] = -0.6290699958801254
$arr[17,1] = 1.975620031356814
$arr[17,2] = -0.5599067687988319
$arr[17,3] = -0.0335975885391235
$arr[17,4] = -0.0216857157647609
$arr[17,5] = -0.110566608607769
$arr[18,0] = -1.261460304260257
$arr[18,1] = 1.682102203369138
$arr[18,2] = -0.007327961921688053
$arr[18,3] = 0.0045814891345798
$arr[18,4] = 0.0555887371301651
$arr[18,5] = 0.012980886735022
$arr[19,0] = -0.9165861129760703
$arr[19,1] = 1.882461953163149
$arr[19,2] = 0.1453821629285811
$arr[19,3] = -0.0044287731871008
$arr[19,4] = 0.0375682115554809
$arr[19,5] = 0.0296269636601209

$ws.Range("C2:H21").Value = $arr

# Remove the trailing row that is no longer needed after the shift
$ws.Range("A22:H22").ClearContents()
